$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update B2/B3 numeric values ---
$ws.Range("B2").Value = 113
$ws.Range("B3").Value = 114

# --- Row 9: new number (right-aligned style like B2/B3) + new JSON string in C9 ---
$ws.Range("B2").Copy()
$ws.Range("B9").PasteSpecial(-4122)
$ws.Range("B9").Value = 115

$json1 = "{`n    ""username"": ""user676767"",`n    ""fromplace"": ""TVM"",`n    ""toplace"": ""BGRLR"",`n    ""email"": ""user676767@gmail.com"",`n    ""price"": 1245.0,`n    ""traveldate"": null,`n    ""pincode"": ""321433""`n}`n`n"
$ws.Range("C9").Value = $json1
$ws.Rows.Item(9).RowHeight = 14.25

# --- Row 10: new number (right-aligned style like B2/B3) + new JSON string in C10 ---
$ws.Range("B2").Copy()
$ws.Range("B10").PasteSpecial(-4122)
$ws.Range("B10").Value = 116

$json2 = "{`n    ""username"": ""user676767"",`n    ""fromplace"": ""Chennai"",`n    ""toplace"": ""Cochin"",`n    ""email"": ""user676767@gmail.com"",`n    ""price"": 6734.0,`n    ""traveldate"": null,`n    ""pincode"": ""321456""`n}`n`n"
$ws.Range("C10").Value = $json2
$ws.Rows.Item(10).RowHeight = 14.25

# --- B12/B13 numeric values ---
$ws.Range("B12").Value = 117
$ws.Range("B13").Value = 118

# --- Update the active selection to C14 ---
$ws.Range("C14").Select()
